# Add new negative-manifest rows (n21-n24) to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("n21", "n21_e68_321_2_1.jpeg", "True", "no_meltpatch", "negative"),
    @("n22", "n22_e67_321_2_2.jpeg", "True", "no_meltpatch", "negative"),
    @("n23", "n23_e70_321_3_1.jpeg", "True", "no_meltpatch", "negative"),
    @("n24", "n24_e65_321_1_2.jpeg", "True", "no_meltpatch", "negative")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        # Prefix with a leading apostrophe so values like "True" are stored
        # as literal text rather than being auto-coerced into a Boolean.
        $cell.Value = "'" + $rowData[$c]
        # Clear the "quote prefix" formatting the apostrophe leaves behind
        # so the new cells keep the same (default) style as existing ones.
        $cell.Style = "Normal"
    }
}

$wb.Save()
